$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append academic distinctions to the PhD (row 2) and MSc (row 4) entries.
$ws.Range("A2").Value = "PhD in Neuroscience  - \textbf{\textit{Summa Cum Laude}}"
$ws.Range("A4").Value = "Psychological Research Methods (Evolutionary Psychology) MSc  - \textbf{\textit{Trabajo de grado meritorio}}"

# Move the active selection to A4, matching the saved view state.
$ws.Range("A4").Select()
